$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.855.55'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.26%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.869.65'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.36%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7354'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.31%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '241.85'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.18%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.000'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3153'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.38%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '24.73'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.57%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07091'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.47%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08383'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -8.53%  '

$ws.Range("E12").Value = '  -3.06%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.413'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.69%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.841.87'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.50%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.55'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.43%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.851.24'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.39%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.046'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.01%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.57'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.09%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '242.95'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.61%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007832'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.93%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9997'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.08%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.116.14'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.94%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.903'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.12%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1567'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.54%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.310'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.55%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '164.17'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.80%  '

$ws.Range("E28").Value = '  -1.45%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.015'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.66%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.477'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.70%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.623'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.94%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.531'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.43%  '

$ws.Range("E33").Value = '  +4.12%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05332'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.44%  '

$ws.Range("E35").Value = '  -1.31%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7533'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.45%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.000'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.24%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.701'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.45%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01952'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.55%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.752'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.41%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4466'
$ws.Range("D41").Style = "Normal"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.106.80'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.16%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.082'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.43%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '72.14'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.79%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8606'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.40%  '

$ws.Range("E46").Value = '  +0.10%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '103.01'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.17%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.703'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.83%  '

$ws.Range("E49").Value = '  -2.79%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.059'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.68%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.014.70'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.29%  '

